$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook is a localization-status report with 3 sheets:
#   Overview (sheet1) - one row per handed-off source file
#   zh-cn    (sheet2) - per-file status for the zh-cn target
#   de-de    (sheet3) - per-file status for the de-de target
# A new file (66d02b38-e0f4-4437-a186-3f0a73296a7a.md) was handed off, so a
# new row (row 3) needs to be appended to every sheet, together with its
# hyperlink and the corresponding table/autofilter range growing by one row.
# ---------------------------------------------------------------------------

$commit = "652e3c1f3f099dfd9354b6e9f2a24f741c200bd2"
$fileGuid = "66d02b38-e0f4-4437-a186-3f0a73296a7a"
$fileName = "$fileGuid.md"
$pathName = "e2e\$fileGuid.md"
$url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit/e2e/$fileName"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Write the new row's cell values left-to-right, sheet-by-sheet (Overview,
# then zh-cn, then de-de) so freshly introduced strings are interned by the
# engine in the same left-to-right / sheet-by-sheet order the original
# report generator produced them in.
# ---------------------------------------------------------------------------

# Overview sheet, row 3
$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("B3").Value = $pathName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-17 06:35:13"

# zh-cn sheet, row 3
$wsZhCn.Range("A3").Value = $fileName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "$fileGuid.855381a1591b842cd73685f1b35ac7d23105dd09.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-17 06:35:08"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

# de-de sheet, row 3
$wsDeDe.Range("A3").Value = $fileName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "$fileGuid.855381a1591b842cd73685f1b35ac7d23105dd09.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-17 06:35:13"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

# ---------------------------------------------------------------------------
# Formatting: reuse the formats from the row above (hyperlink style, date
# display format) so the new row visually matches the existing rows.
# ---------------------------------------------------------------------------
$wsOverview.Range("G3").NumberFormat = $wsOverview.Range("G2").NumberFormat
$wsZhCn.Range("H3").NumberFormat = $wsZhCn.Range("H2").NumberFormat
$wsZhCn.Range("K3").NumberFormat = $wsZhCn.Range("K2").NumberFormat
$wsDeDe.Range("H3").NumberFormat = $wsDeDe.Range("H2").NumberFormat
$wsDeDe.Range("K3").NumberFormat = $wsDeDe.Range("K2").NumberFormat

# ---------------------------------------------------------------------------
# Hyperlinks for the new handed-off file.
# ---------------------------------------------------------------------------
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $url, "", "", $pathName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $url, "", "", $fileName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $url, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------------
# Grow each table / autofilter range to include the new row.
# ---------------------------------------------------------------------------
$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G3"))
$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P3"))
$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P3"))
